$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.035.16"
$ws.Range("E2").Value = "  -3.60%  "
$ws.Range("D3").Value = "3.437.06"
$ws.Range("E3").Value = "  +2.45%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'236.16"
$ws.Range("E5").Value = "  -6.39%  "
$ws.Range("D6").Value = "'637.00"
$ws.Range("E6").Value = "  -3.45%  "
$ws.Range("D7").Value = "'1.43"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("D8").Value = "'0.397"
$ws.Range("E8").Value = "  -7.32%  "
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("D10").Value = "'0.965"
$ws.Range("E10").Value = "  -6.12%  "
$ws.Range("D11").Value = "3.437.69"
$ws.Range("E11").Value = "  +2.56%  "
$ws.Range("D12").Value = "'42.08"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("E13").Value = "  -5.20%  "
$ws.Range("D14").Value = "'6.15"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "93.888.59"
$ws.Range("E15").Value = "  -3.45%  "
$ws.Range("D16").Value = "4.085.26"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("E17").Value = "  -1.95%  "
$ws.Range("D18").Value = "'8.36"
$ws.Range("E18").Value = "  -6.30%  "
$ws.Range("D19").Value = "3.456.41"
$ws.Range("E19").Value = "  +3.66%  "
$ws.Range("D20").Value = "'17.60"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").Value = "'11.26"
$ws.Range("E21").Value = "  +3.98%  "
$ws.Range("E22").Value = "  -11.34%  "
$ws.Range("D23").Value = "'497.71"
$ws.Range("E23").Value = "  -2.94%  "
$ws.Range("D24").Value = "'3.17"
$ws.Range("E24").Value = "  -6.19%  "
$ws.Range("D25").Value = "'6.57"
$ws.Range("E25").Value = "  -1.49%  "
$ws.Range("E26").Value = "  -6.69%  "
$ws.Range("D27").Value = "'90.87"
$ws.Range("E27").Value = "  -6.88%  "
$ws.Range("D28").Value = "3.619.74"
$ws.Range("E28").Value = "  +2.49%  "
$ws.Range("D29").Value = "'11.95"
$ws.Range("E29").Value = "  -2.68%  "
$ws.Range("D30").Value = "'11.74"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").Value = "'2.73"
$ws.Range("E32").Value = "  +6.82%  "
$ws.Range("D33").Value = "'0.136"
$ws.Range("E33").Value = "  -8.92%  "
$ws.Range("D34").Value = "'0.182"
$ws.Range("E34").Value = "  -4.79%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").Value = "'30.17"
$ws.Range("E36").Value = "  +5.27%  "
$ws.Range("D37").Value = "'0.557"
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("D38").Value = "'546.24"
$ws.Range("E38").Value = "  +6.04%  "
$ws.Range("D39").Value = "'7.62"
$ws.Range("E39").Value = "  -4.77%  "
$ws.Range("E40").Value = "  -5.18%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'0.936"
$ws.Range("E41").Value = "  +10.21%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").Value = "'24.05"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").Value = "'1.68"
$ws.Range("E45").Value = "  -1.39%  "
$ws.Range("E46").Value = "  -5.92%  "
$ws.Range("D47").Value = "'5.52"
$ws.Range("E47").Value = "  -3.36%  "
$ws.Range("D48").Value = "'3.48"
$ws.Range("E48").Value = "  -4.91%  "
$ws.Range("E49").Value = "  +6.10%  "
$ws.Range("D50").Value = "'53.09"
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("E51").Value = "  +0.10%  "
